$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -12.621
$ws.Range("E2").Value = 16.597
$ws.Range("A3").Value = -21.837
$ws.Range("D3").Value = -7.561999999999999
$ws.Range("E6").Value = 16.837
$ws.Range("D12").Value = -7.611
$ws.Range("A14").Value = -21.913
$ws.Range("A16").Value = -22.075
$ws.Range("C18").Value = -12.282
$ws.Range("E19").Value = 16.434
$ws.Range("A21").Value = -20.217
$ws.Range("A23").Value = -20.299
$ws.Range("C24").Value = -11.891
$ws.Range("D24").Value = -7.526999999999999
$ws.Range("E24").Value = 16.952
$ws.Range("A25").Value = -22.037
$ws.Range("C25").Value = -13.549
$ws.Range("D25").Value = -8.716999999999999
$ws.Range("A26").Value = -21.435
$ws.Range("C27").Value = -13.128
$ws.Range("E27").Value = 16.658
$ws.Range("A29").Value = -21.219
$ws.Range("C30").Value = -12.686
$ws.Range("E30").Value = 16.887
$ws.Range("C31").Value = -12.126
$ws.Range("E31").Value = 17.123
$ws.Range("E33").Value = 17.573
$ws.Range("C39").Value = -12.702
$ws.Range("A40").Value = -19.938
$ws.Range("D41").Value = -8.195
$ws.Range("C42").Value = -12.852
$ws.Range("E42").Value = 16.606
$ws.Range("C48").Value = -11.598
$ws.Range("D50").Value = -8.233000000000001
$ws.Range("C51").Value = -11.564
$ws.Range("C52").Value = -11.606
$ws.Range("A53").Value = -21.924
$ws.Range("D53").Value = -7.815
$ws.Range("C55").Value = -13.341
$ws.Range("E55").Value = 16.057
$ws.Range("C56").Value = -12.362
$ws.Range("D56").Value = -8.358000000000001
$ws.Range("A57").Value = -21.891
$ws.Range("C57").Value = -12.926
$ws.Range("D57").Value = -8.822999999999999
$ws.Range("D58").Value = -8.318
$ws.Range("E58").Value = 16.589
$ws.Range("A59").Value = -22.323
$ws.Range("C60").Value = -12.296
$ws.Range("D61").Value = -7.708
$ws.Range("D63").Value = -7.885999999999998
$ws.Range("D64").Value = -7.930000000000001
$ws.Range("A65").Value = -21.421
$ws.Range("E65").Value = 17.114
$ws.Range("A69").Value = -21.833
$ws.Range("D70").Value = -7.598000000000001
$ws.Range("E70").Value = 17.612
$ws.Range("D72").Value = -7.356999999999999
$ws.Range("C73").Value = -12.995
$ws.Range("C74").Value = -12.077
$ws.Range("E74").Value = 16.575
$ws.Range("E75").Value = 16.76
$ws.Range("A79").Value = -20.849
$ws.Range("A83").Value = -21.938
$ws.Range("E83").Value = 16.874
$ws.Range("E84").Value = 17.027
$ws.Range("D86").Value = -8.367000000000001
$ws.Range("E86").Value = 16.828
$ws.Range("C89").Value = -12.305
$ws.Range("D89").Value = -7.857999999999999
$ws.Range("C90").Value = -13.047
$ws.Range("A91").Value = -21.508
$ws.Range("C92").Value = -11.915
$ws.Range("A93").Value = -21.203
$ws.Range("E96").Value = 16.876
$ws.Range("E97").Value = 16.846
$ws.Range("D98").Value = -8.612
$ws.Range("A100").Value = -22.017
$ws.Range("D100").Value = -8.690999999999999
$ws.Range("D102").Value = -7.859
